$wb = $excel.ActiveWorkbook

$wsPersonnel = $wb.Worksheets.Item("Personnel")
$wsKeywords  = $wb.Worksheets.Item("Keywords")

# Add missing ORCID userId for Katherine Qi (row 4, column F = userId)
$wsPersonnel.Activate()
$wsPersonnel.Range("F4").Value = "0000-0002-6839-2579"
$wsPersonnel.Range("F5").Select()

# Update the keyword thesaurus title text
$wsKeywords.Activate()
$wsKeywords.Range("B11").Value = "OBO Open Biological and Biomedical Ontology"
$wsKeywords.Range("A11").Select()
